# "Generate Report for Handback"
#
# This localization-status report is regenerated: the zh-cn and de-de
# localization targets have been handed back and are now in sync with
# en-US, so the report is refreshed accordingly:
#   - Status cells move from "Ready for handoff" to
#     "Handed back: in sync with en-US" (this text is shared by the
#     Overview sheet's zh-cn/de-de status cells and by the Status column
#     on each per-language sheet).
#   - The "Latest Handback DateTime" for each language is refreshed to
#     the handback timestamp.
#   - The stale "version mismatch" Error Detail is cleared now that the
#     handback files are in sync.
#   - The Status / Error Detail columns are resized to fit the new
#     report contents.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status cells ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-15 16:46:28"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-15 16:46:36"
$wsDeDe.Range("P2").Value = ""

# --- Column width refresh (Status column widened, Error Detail narrowed) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
